# Updated cryptos list on Tue Sep  5 18:54:29 UTC 2023 with GitHub Actions
#
# Note: several values in column D ("Price") are digit strings that look like
# plain numbers to Excel's automatic type inference (e.g. "1.002", "215.51").
# The source data must stay as TEXT (as in the original workbook, which stores
# every Price/Volume cell as an inline string), so those are written with a
# leading apostrophe to force text entry; the apostrophe itself is not stored
# in the resulting value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "25.795.35"
$ws.Cells.Item(2, 5).Value = "  -0.58%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "1.638.54"
$ws.Cells.Item(3, 5).Value = "  +0.30%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 5).Value = "  +0.11%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "'215.51"
$ws.Cells.Item(5, 5).Value = "  -0.23%  "

# Row 6 - XRP
$ws.Cells.Item(6, 4).Value = "'0.5062"
$ws.Cells.Item(6, 5).Value = "  -1.08%  "

# Row 7 - USDC
$ws.Cells.Item(7, 5).Value = "  +0.06%  "

# Row 8 - Cardano
$ws.Cells.Item(8, 4).Value = "'0.2582"
$ws.Cells.Item(8, 5).Value = "  +0.11%  "

# Row 9 - Dogecoin
$ws.Cells.Item(9, 4).Value = "'0.06427"

# Row 10 - Solana
$ws.Cells.Item(10, 4).Value = "'20.38"
$ws.Cells.Item(10, 5).Value = "  +4.41%  "

# Row 11 - TRON
$ws.Cells.Item(11, 4).Value = "'0.07789"
$ws.Cells.Item(11, 5).Value = "  +0.07%  "

# Row 12 - Polkadot
$ws.Cells.Item(12, 4).Value = "'4.268"
$ws.Cells.Item(12, 5).Value = "  -0.36%  "

# Row 13 & 14 - swapped: WrappedliquidstakedEther2.0 <-> WrappedEther
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.638.21"
$ws.Cells.Item(13, 5).Value = "  +0.23%  "

$ws.Cells.Item(14, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(14, 4).Value = "1.863.83"
$ws.Cells.Item(14, 5).Value = "  +0.26%  "

# Row 15 - Polygon
$ws.Cells.Item(15, 4).Value = "'0.5613"
$ws.Cells.Item(15, 5).Value = "  +1.62%  "

# Row 16 - ShibaInu
$ws.Cells.Item(16, 4).Value = "0.0₅7651"
$ws.Cells.Item(16, 5).Value = "  -0.05%  "

# Row 17 - Litecoin
$ws.Cells.Item(17, 4).Value = "'63.31"
$ws.Cells.Item(17, 5).Value = "  -1.05%  "

# Row 18 - WrappedBTC
$ws.Cells.Item(18, 4).Value = "25.816.22"
$ws.Cells.Item(18, 5).Value = "  -0.57%  "

# Row 19 - Dai
$ws.Cells.Item(19, 5).Value = "  -0.04%  "

# Row 20 & 21 - swapped: BitcoinCash <-> Uniswap
$ws.Cells.Item(20, 2).Value = "Uniswap"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(20, 4).Value = "'4.380"
$ws.Cells.Item(20, 5).Value = "  -1.02%  "

$ws.Cells.Item(21, 2).Value = "BitcoinCash"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(21, 4).Value = "'192.89"
$ws.Cells.Item(21, 5).Value = "  -1.16%  "

# Row 22 - Avalanche
$ws.Cells.Item(22, 4).Value = "'9.924"
$ws.Cells.Item(22, 5).Value = "  +0.52%  "

# Row 23 - Chainlink
$ws.Cells.Item(23, 4).Value = "'6.138"
$ws.Cells.Item(23, 5).Value = "  +1.48%  "

# Row 24 - BinanceUSD
$ws.Cells.Item(24, 4).Value = "'1.001"
$ws.Cells.Item(24, 5).Value = "  -0.02%  "

# Row 25 - Toncoin
$ws.Cells.Item(25, 4).Value = "'1.793"
$ws.Cells.Item(25, 5).Value = "  -5.21%  "

# Row 26 - Monero
$ws.Cells.Item(26, 4).Value = "'140.85"
$ws.Cells.Item(26, 5).Value = "  -0.69%  "

# Row 27 - Stellar
$ws.Cells.Item(27, 4).Value = "'0.1237"
$ws.Cells.Item(27, 5).Value = "  -1.83%  "

# Row 28 - Cosmos
$ws.Cells.Item(28, 4).Value = "'6.821"
$ws.Cells.Item(28, 5).Value = "  +0.94%  "

# Row 29 - EthereumClassic
$ws.Cells.Item(29, 5).Value = "  -0.37%  "

# Row 30 - PancakeSwap
$ws.Cells.Item(30, 4).Value = "'1.244"
$ws.Cells.Item(30, 5).Value = "  +0.21%  "

# Row 31 - Hedera
$ws.Cells.Item(31, 5).Value = "  +1.32%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Cells.Item(32, 4).Value = "'3.281"
$ws.Cells.Item(32, 5).Value = "  +1.09%  "

# Row 33 - Filecoin
$ws.Cells.Item(33, 4).Value = "'3.236"

# Row 34 - LidoDAOToken
$ws.Cells.Item(34, 4).Value = "'1.570"

# Row 35 - HuobiToken
$ws.Cells.Item(35, 4).Value = "'2.382"
$ws.Cells.Item(35, 5).Value = "  +0.44%  "

# Row 36 - ARBITRUM
$ws.Cells.Item(36, 4).Value = "'0.9037"
$ws.Cells.Item(36, 5).Value = "  +0.59%  "

# Row 37 & 38 - swapped: ImmutableX <-> MXToken
$ws.Cells.Item(37, 2).Value = "MXToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(37, 4).Value = "'2.571"
$ws.Cells.Item(37, 5).Value = "  +1.22%  "

$ws.Cells.Item(38, 2).Value = "ImmutableX"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(38, 4).Value = "'0.5559"
$ws.Cells.Item(38, 5).Value = "  +0.73%  "

# Row 39 - Maker
$ws.Cells.Item(39, 4).Value = "1.131.87"
$ws.Cells.Item(39, 5).Value = "  +1.18%  "

# Row 40 - VeChain
$ws.Cells.Item(40, 5).Value = "  +0.60%  "

# Row 41 - PaxDollar
$ws.Cells.Item(41, 4).Value = "'0.9951"
$ws.Cells.Item(41, 5).Value = "  -0.50%  "

# Row 42 - FraxShare
$ws.Cells.Item(42, 4).Value = "'5.481"
$ws.Cells.Item(42, 5).Value = "  -2.02%  "

# Row 43 - TrustWalletToken
$ws.Cells.Item(43, 5).Value = "  +0.85%  "

# Row 44 - Quant
$ws.Cells.Item(44, 4).Value = "'98.93"
$ws.Cells.Item(44, 5).Value = "  +1.39%  "

# Row 45 - RocketPoolETH
$ws.Cells.Item(45, 4).Value = "1.773.92"
$ws.Cells.Item(45, 5).Value = "  +0.22%  "

# Row 46 - BabyDogeCoin
$ws.Cells.Item(46, 4).Value = "0.0₈109"
$ws.Cells.Item(46, 5).Value = "  -6.68%  "

# Row 47 - Aave
$ws.Cells.Item(47, 4).Value = "'55.63"
$ws.Cells.Item(47, 5).Value = "  +1.62%  "

# Row 48 - Mantle
$ws.Cells.Item(48, 4).Value = "'0.4273"
$ws.Cells.Item(48, 5).Value = "  -3.88%  "

# Row 49 - EnergySwap
$ws.Cells.Item(49, 4).Value = "'7.776"
$ws.Cells.Item(49, 5).Value = "  +2.75%  "

# Row 50 - Cronos
$ws.Cells.Item(50, 4).Value = "'0.05033"
$ws.Cells.Item(50, 5).Value = "  -1.90%  "

# Row 51 - Frax
$ws.Cells.Item(51, 4).Value = "'0.9983"
$ws.Cells.Item(51, 5).Value = "  -0.45%  "
